$d = $word.ActiveDocument

# Insert six new list paragraphs right after paragraph 5
# ("Create Days Objects for each of the days selected ...").
$anchor = $d.Paragraphs(5)
$anchor.Range.InsertParagraphAfter()

$p1 = $d.Paragraphs(6)
$p1.Range.Text = "Each Day provides the framework for the Lemonade stand simulation game."

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(7)
$p2.Range.Text = "First set correct number of Days."
$p2.Range.ListFormat.ListLevelNumber = 3

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(8)
$p3.Range.Text = "Randomly select a temperature for each day."
$p3.Range.ListFormat.ListLevelNumber = 3

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(9)
$p4.Range.Text = "Randomly select a weather condition."
$p4.Range.ListFormat.ListLevelNumber = 3

$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(10)
$p5.Range.Text = "Instantiate Days with unique weather condition and temperature"
$p5.Range.ListFormat.ListLevelNumber = 3

$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs(11)
$p6.Range.Text = "Based on temperature and weather condition create a number of customers."
$p6.Range.ListFormat.ListLevelNumber = 3

Write-Output $d.Paragraphs.Count
